$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.071.13'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.563.41'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.54'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.491'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.93'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.59%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.784.84'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.565.35'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.77'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.057.14'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.97'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0701'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.49'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.84'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.03'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("E28").Value = '  +1.42%  '
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.34%  '
$ws.Range("E31").Value = '  +0.42%  '
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.434.61'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("E35").Value = '  +2.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.81'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.84%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.40'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.699.37'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0518'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0958'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.23%  '
